# Add a new acronym-key row for "IiCFfNPdtTI" (Increase in Capacity Factors for
# New Plants due to Technological Improvements) into the "elec" group of the
# "Key to Variables" sheet, right above the "MPCbS" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new blank row at row 93 (pushes MPCbS and everything below down by one).
$ws.Rows.Item(93).Insert()

# Copy the formatting (fill/wrap etc.) of the surrounding table so the new row
# matches the rest of the "elec" block before we fill in its own values.
$ws.Range("A92:C92").Copy()
$ws.Range("A93:C93").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F95").Copy()
$ws.Range("F93").PasteSpecial(-4122)       # xlPasteFormats (matches the "low" color used elsewhere)

$excel.CutCopyMode = 0

# Fill in the new row's values.
$ws.Range("A93").Value2 = "elec"
$ws.Range("B93").Value2 = "IiCFfNPdtTI"
$ws.Range("C93").Value2 = "Increase in Capacity Factors for New Plants due to Technological Improvements"
$ws.Range("F93").Value2 = "low"

# The meaning text wraps onto two lines at the sheet's column width, so give the
# row the extra height (matches every other two-line row in the sheet).
$ws.Rows.Item(93).RowHeight = 30

# Leave the view selected on the row we just added.
$ws.Activate()
$ws.Range("A93").Select()
